# "First long run OK"
# Update the two existing rows (dates / amounts) and append the remaining
# campaign rows to the tracking sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 2 - "Personas con Discapacidad": shift dates, widen the +2 -> +4
# window, and bump the requested volume.
# ---------------------------------------------------------------------
$ws.Range("B2").Value = 44166
$ws.Range("C2").Formula = "=B2+4"
$ws.Range("D2").Value = 40000

# ---------------------------------------------------------------------
# Row 3 - "#DiaInternacionalAsperger": same treatment.
# ---------------------------------------------------------------------
$ws.Range("B3").Value = 43877
$ws.Range("C3").Formula = "=B3+4"
$ws.Range("D3").Value = 40000

# ---------------------------------------------------------------------
# New rows 4-11. Clone the date formatting from the existing B2/C2 cells
# (style carries the numFmtId=14 date format) before dropping in values,
# then enter the query names in the specific order that matches how they
# were first typed (keeps the shared-string table order stable).
# ---------------------------------------------------------------------
$ws.Range("B2").Copy()
$ws.Range("B4:B11").PasteSpecial(-4122)

$ws.Range("C2").Copy()
$ws.Range("C4:C11").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Range("A6").Value = "#DíaMundialdelaVisión"
$ws.Range("A8").Value = "#DiaMundialAutismo"
$ws.Range("A10").Value = "#DiaMundialEnfermedadesRaras"
$ws.Range("A5").Value = "#UnFuturoMejor"
$ws.Range("A4").Value = "#TDAH"
$ws.Range("A7").Value = "#lacausaquenosune"
$ws.Range("A9").Value = "#calcetinesdesparejados"
$ws.Range("A11").Value = "#DiaNacionalEB"

$ws.Range("B4").Value = 44129
$ws.Range("C4").Formula = "=B4+4"
$ws.Range("D4").Value = 40000

$ws.Range("B5").Value = 44109
$ws.Range("C5").Formula = "=B5+4"
$ws.Range("D5").Value = 40000

$ws.Range("B6").Value = 44110
$ws.Range("C6").Formula = "=B6+4"
$ws.Range("D6").Value = 40000

$ws.Range("B7").Value = 43734
$ws.Range("C7").Formula = "=B7+4"
$ws.Range("D7").Value = 40000

$ws.Range("B8").Value = 43920
$ws.Range("C8").Formula = "=B8+4"
$ws.Range("D8").Value = 40000

$ws.Range("B9").Value = 43909
$ws.Range("C9").Formula = "=B9+4"
$ws.Range("D9").Value = 40000

$ws.Range("B10").Value = 43888
$ws.Range("C10").Formula = "=B10+4"
$ws.Range("D10").Value = 40000

$ws.Range("B11").Value = 43788
$ws.Range("C11").Formula = "=B11+4"
$ws.Range("D11").Value = 40000

# Rows 6 & 7 hold the longer hashtags; Excel wrapped those two labels and
# grew the row height to fit the extra line.
$ws.Range("A6").WrapText = $true
$ws.Range("A7").WrapText = $true
$ws.Rows.Item(6).RowHeight = 17
$ws.Rows.Item(7).RowHeight = 17

# Column sizing: column A widened to fit the longer labels, column B
# picked up an explicit best-fit width once the date column got busy.
$ws.Columns.Item(1).ColumnWidth = 29.5
$ws.Columns.Item(2).ColumnWidth = 8.8

# Final selection the author left the sheet on.
[void]$ws.Range("D19").Select()
